$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 "Save", matching style of G1 (bold header with border)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# New data cells H2 and H3 with value 0
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0

Write-Host "Done"
